$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "56+7="
$t.Cell(1,2).Range.Text = "97-19="
$t.Cell(1,3).Range.Text = "92-20="
$t.Cell(1,4).Range.Text = "18+5="
$t.Cell(1,5).Range.Text = "24+74="
$t.Cell(2,1).Range.Text = "99-34="
$t.Cell(2,2).Range.Text = "89-63="
$t.Cell(2,3).Range.Text = "8+8="
$t.Cell(2,4).Range.Text = "9+3="
$t.Cell(2,5).Range.Text = "74-8="
$t.Cell(3,1).Range.Text = "81-23="
$t.Cell(3,2).Range.Text = "22+0="
$t.Cell(3,3).Range.Text = "43-33="
$t.Cell(3,4).Range.Text = "96-69="
$t.Cell(3,5).Range.Text = "61+23="
$t.Cell(4,1).Range.Text = "32+53="
$t.Cell(4,2).Range.Text = "49-43="
$t.Cell(4,3).Range.Text = "86-40="
$t.Cell(4,4).Range.Text = "88-6="
$t.Cell(4,5).Range.Text = "99-73="
$t.Cell(5,1).Range.Text = "72-29="
$t.Cell(5,2).Range.Text = "67-21="
$t.Cell(5,3).Range.Text = "89-84="
$t.Cell(5,4).Range.Text = "39+20="
$t.Cell(5,5).Range.Text = "87-72="
$t.Cell(6,1).Range.Text = "69-65="
$t.Cell(6,2).Range.Text = "46+0="
$t.Cell(6,3).Range.Text = "17+5="
$t.Cell(6,4).Range.Text = "53-35="
$t.Cell(6,5).Range.Text = "0+77="
$t.Cell(7,1).Range.Text = "21+73="
$t.Cell(7,2).Range.Text = "96+2="
$t.Cell(7,3).Range.Text = "32+44="
$t.Cell(7,4).Range.Text = "17+4="
$t.Cell(7,5).Range.Text = "7-3="
$t.Cell(8,1).Range.Text = "28+53="
$t.Cell(8,2).Range.Text = "64+10="
$t.Cell(8,3).Range.Text = "47-5="
$t.Cell(8,4).Range.Text = "12+1="
$t.Cell(8,5).Range.Text = "70-29="
$t.Cell(9,1).Range.Text = "20+46="
$t.Cell(9,2).Range.Text = "14+25="
$t.Cell(9,3).Range.Text = "37+0="
$t.Cell(9,4).Range.Text = "67+22="
$t.Cell(9,5).Range.Text = "59+10="
$t.Cell(10,1).Range.Text = "54-48="
$t.Cell(10,2).Range.Text = "84-20="
$t.Cell(10,3).Range.Text = "85-72="
$t.Cell(10,4).Range.Text = "96-87="
$t.Cell(10,5).Range.Text = "66-49="
$t.Cell(11,1).Range.Text = "42+53="
$t.Cell(11,2).Range.Text = "80-65="
$t.Cell(11,3).Range.Text = "18-2="
$t.Cell(11,4).Range.Text = "52+31="
$t.Cell(11,5).Range.Text = "16+59="
$t.Cell(12,1).Range.Text = "39-22="
$t.Cell(12,2).Range.Text = "67-19="
$t.Cell(12,3).Range.Text = "97-31="
$t.Cell(12,4).Range.Text = "68-34="
$t.Cell(12,5).Range.Text = "6+73="
$t.Cell(13,1).Range.Text = "49-18="
$t.Cell(13,2).Range.Text = "17+37="
$t.Cell(13,3).Range.Text = "56-12="
$t.Cell(13,4).Range.Text = "96-60="
$t.Cell(13,5).Range.Text = "86-6="
$t.Cell(14,1).Range.Text = "9+59="
$t.Cell(14,2).Range.Text = "52+13="
$t.Cell(14,3).Range.Text = "76-1="
$t.Cell(14,4).Range.Text = "69-18="
$t.Cell(14,5).Range.Text = "91-74="
$t.Cell(15,1).Range.Text = "32+58="
$t.Cell(15,2).Range.Text = "57-22="
$t.Cell(15,3).Range.Text = "68+27="
$t.Cell(15,4).Range.Text = "37+57="
$t.Cell(15,5).Range.Text = "80-14="
$t.Cell(16,1).Range.Text = "94-0="
$t.Cell(16,2).Range.Text = "64-56="
$t.Cell(16,3).Range.Text = "50+27="
$t.Cell(16,4).Range.Text = "71-48="
$t.Cell(16,5).Range.Text = "17+70="
$t.Cell(17,1).Range.Text = "73-10="
$t.Cell(17,2).Range.Text = "35-19="
$t.Cell(17,3).Range.Text = "18+30="
$t.Cell(17,4).Range.Text = "72-4="
$t.Cell(17,5).Range.Text = "12+79="
$t.Cell(18,1).Range.Text = "54+23="
$t.Cell(18,2).Range.Text = "54+35="
$t.Cell(18,3).Range.Text = "44+17="
$t.Cell(18,4).Range.Text = "4+82="
$t.Cell(18,5).Range.Text = "21+18="
$t.Cell(19,1).Range.Text = "45-21="
$t.Cell(19,2).Range.Text = "69-26="
$t.Cell(19,3).Range.Text = "39+32="
$t.Cell(19,4).Range.Text = "71+23="
$t.Cell(19,5).Range.Text = "49-49="
$t.Cell(20,1).Range.Text = "58+15="
$t.Cell(20,2).Range.Text = "8+8="
$t.Cell(20,3).Range.Text = "30+55="
$t.Cell(20,4).Range.Text = "26+70="
$t.Cell(20,5).Range.Text = "96-33="
